# Apply the LinuxForHealth re-brand / version bump edits to the
# StructureDefinition-communication-template workbook.
#
# Sheet "Metadata" (A/B key-value table):
#   URL        -> http://linuxforhealth.org/fhir/cdm/StructureDefinition/communication-template
#   Version    -> 8.0.0
#   Date       -> 2022-11-10T16:00:46+00:00
#   Publisher  -> LinuxForHealth Team
#
# Sheet "Elements": the root "Extension" row's Constraint(s) cell (AI2) is
# cleared - it was a stray duplicate of the ele-1/ext-1 constraint text that
# correctly lives on the "Extension.extension" row (AI4) instead.

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/communication-template"
$metadata.Range("B3").Value = "8.0.0"
$metadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$metadata.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""
